$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: B1 changes from "Actual" to "Value"
$ws.Range("B1").Value = "Value"

# Restore the selection to C1 (was C7 before the edit)
$ws.Range("C1").Select()
